# Update "想去人数" (F) and "最低票价" (G) figures on the 展览 and 全部类型
# sheets to reflect the latest scrape output.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 2136
    $ws.Range("F3").Value = 619
    $ws.Range("F4").Value = 1516
    $ws.Range("F5").Value = 7204
    # Keep "58" as text (matches the existing text-stored ticket prices in
    # column G); a bare numeric-looking string would otherwise be
    # auto-coerced to a number, so prefix it like Excel's quote-prefix entry.
    $ws.Range("G5").Value = "'58"
    $ws.Range("F6").Value = 177
    $ws.Range("F7").Value = 150
}
